# address point about linearity of filters
#
# The reviewer-response paragraph that used to read:
#   "Seelig and Jayaraman verified the RFs with the standard white noise
#    stimuli technique (their refs 37,41). This is now mentioned in the
#    text."
# (shown with a yellow highlight, using the "TextBody" paragraph style)
# is replaced with the authors' fuller answer about linear filters, no
# longer highlighted and using the surrounding "PreformattedText" style
# that the rest of the reply uses.

$d = $word.ActiveDocument

$oldSeelig = "Seelig and Jayaraman verified the RFs with the standard white noise stimuli technique (their refs 37,41). This is now mentioned in the text."
$newSeelig = "Seelig and Jayaraman used the RFs in this way as linear convolutions (to verify the RF forms). We have added some discussion of this to Materials and methods. We have additionally justified the use of linear filters in the main text (in the Results section)."

$d.Content.Find.Execute($oldSeelig, $false, $false, $false, $false, $false, `
                         $true, 1, $false, $newSeelig, 2) | Out-Null

# Re-locate the paragraph (now holding the new wording) and bring its
# paragraph style / run formatting in line with the rest of the reply:
# drop the yellow highlight, drop the explicit "Liberation Mono" run
# font (the PreformattedText style already supplies it) and switch the
# paragraph from "TextBody" to "PreformattedText".
$rng = $d.Content
$rng.Find.Execute($newSeelig, $false, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0) | Out-Null

$rng.HighlightColorIndex = 0
$rng.Font.Name = ""
$rng.Paragraphs(1).Style = "PreformattedText"
$rng.Font.Color = 655360
$rng.Font.Size = 10

# Minor punctuation fix elsewhere in the document: "performs better!" ->
# "performs better."
$d.Content.Find.Execute("performs better!", $false, $false, $false, $false, `
                         $false, $true, 1, $false, "performs better.", 2) | Out-Null
